# Updated via Streamlit Approval System
#
# Rows 2-19, columns AK:AO (COST_CENTER, LEDGER_NAME, LEDGER_UNDER, TO, BY)
# are normalized to the text value "0" (they were numeric 0, or blank for
# rows 18-19). Rows 18 and 19 also get their ADJUSTMENT_AMOUNT (AE) filled
# in to match the BASIC_AMOUNT (V) of the pending row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("AK", "AL", "AM", "AN", "AO")

for ($row = 2; $row -le 19; $row++) {
    foreach ($col in $cols) {
        # Leading apostrophe forces Excel to store the value as text "0"
        # rather than the number 0.
        $ws.Range("$col$row").Value = "'0"
    }
}

$ws.Range("AE18").Value = 500
$ws.Range("AE19").Value = 1500
